$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.796.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.486.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.482.85"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.90%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.579"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.01"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.81%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.051.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.45"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "623.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.26%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.821.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.492.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.64%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.878"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.79"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -9.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.33"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.18%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -12.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.34"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -8.95%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.42"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.45%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -9.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "619.22"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.68"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.76%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.11%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.44"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -16.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0442"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.344.04"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.32%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.73%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -11.66%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.10"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.86%  "
